# Fix Employee Middle Name Issue - Fix Employee Import - Other
#
# - Replace the unused "Street2" column (H) with a "City/Suburb" column
#   (header + "Brooklyn" value for the existing employee).
# - Add a new "Gender" column (L) with M/F values.
# - Add a second employee row (Jane Smith) with her own contact info.
# - Turn the Email column entries into mailto: hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace Street2/"Main Street" column with City/Suburb column ---
$ws.Range("H1").Value = "City/Suburb"
$ws.Range("H2").Value = "Brooklyn"

# --- Add the new Gender column after Country, copying the header format ---
$ws.Range("K1").Copy($ws.Range("L1"))
$ws.Range("L1").Value = "Gender"
$ws.Range("L2").Value = "M"
$ws.Range("L3").Value = "F"

# --- Add second employee row (Jane Smith) ---
$ws.Range("A3").Value = "Jane"
$ws.Range("B3").Value = "Smith"
$ws.Range("C3").Value = 9995551213
$ws.Range("D3").Value = 9995551213
$ws.Range("E3").Value = "jane@email.com"
$ws.Range("F3").Value = "janesmith"
$ws.Range("G3").Value = "123 Main Street"
$ws.Range("H3").Value = "Brooklyn"
$ws.Range("I3").Value = "New York"
$ws.Range("J3").Value = 1234
$ws.Range("K3").Value = "United States"

# --- Turn the email addresses into mailto hyperlinks ---
[void]$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:johnsmith@email.com")
[void]$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:jane@email.com")

# --- Update the active selection ---
[void]$ws.Range("H6").Select()
